$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 766351
$ws.Range("E2").Value = 1429263437
$ws.Range("C13").Value = 187874
$ws.Range("E13").Value = 1168847701
$ws.Range("C19").Value = 27520
$ws.Range("E19").Value = 132568759
$ws.Range("C21").Value = 175244
$ws.Range("E21").Value = 316831538
$ws.Range("C41").Value = 126947
$ws.Range("E41").Value = 662714913
$ws.Range("C57").Value = 31599
$ws.Range("E57").Value = 162636444
$ws.Range("C81").Value = 88360
$ws.Range("E81").Value = 499741973
$ws.Range("C88").Value = 71284
$ws.Range("E88").Value = 110329804
$ws.Range("C121").Value = 1306505
$ws.Range("E121").Value = 2275852587
$ws.Range("C129").Value = 633969
$ws.Range("E129").Value = 3437685258
$ws.Range("C132").Value = 586157
$ws.Range("E132").Value = 3476085690
$ws.Range("C151").Value = 39937
$ws.Range("E151").Value = 60395303
$ws.Range("C156").Value = 12418
$ws.Range("E156").Value = 40886248
$ws.Range("C171").Value = 95834
$ws.Range("E171").Value = 490727046
$ws.Range("C178").Value = 515894
$ws.Range("E178").Value = 891232265
$ws.Range("C186").Value = 236848
$ws.Range("E186").Value = 1190253497
$ws.Range("C237").Value = 283332
$ws.Range("E237").Value = 1438703447
$ws.Range("C240").Value = 205945
$ws.Range("E240").Value = 1070200493
